# Rename the "Trait lists" header (column G, row 1) to "Tags".
# All downstream shared-string renumbering and the removal of the now
# unused "Trait lists" string happens automatically on save.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "Tags"

# Move the active selection from H3 to G1, matching the saved view state.
$ws.Range("G1").Select()
